$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns for rows 2-51,
# plus a Coin/Link (B/C) swap between rows 47 and 48 (Arweave <-> Monero).
# Some new Price values parse as plain numbers (e.g. "706.40"); since the sheet
# stores these as text, force text entry via NumberFormat "@" then restore the
# default 'Normal' style so no stray formatting is introduced.

$ws.Range('D2').Value = '71.195.21'
$ws.Range('E2').Value = '  +3.23%  '

$ws.Range('D3').Value = '3.815.67'
$ws.Range('E3').Value = '  +1.19%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '706.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +12.17%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.69'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.01%  '

$ws.Range('D7').Value = '3.813.09'
$ws.Range('E7').Value = '  +1.18%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  +1.29%  '

$ws.Range('E10').Value = '  +3.79%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.45'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.17%  '

$ws.Range('E12').Value = '  +1.22%  '

$ws.Range('E13').Value = '  +9.49%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.67%  '

$ws.Range('D15').Value = '4.457.37'
$ws.Range('E15').Value = '  +1.18%  '

$ws.Range('D16').Value = '3.818.11'
$ws.Range('E16').Value = '  +1.11%  '

$ws.Range('D17').Value = '71.173.73'
$ws.Range('E17').Value = '  +3.21%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.94%  '

$ws.Range('E20').Value = '  +0.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +18.56%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '483.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.46%  '

$ws.Range('E23').Value = '  +2.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.41%  '

$ws.Range('E25').Value = '  +2.84%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.50%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.54%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.99%  '

$ws.Range('D29').Value = '3.967.33'
$ws.Range('E29').Value = '  +1.05%  '

$ws.Range('E30').Value = '  -0.14%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +16.27%  '

$ws.Range('E32').Value = '  +2.23%  '

$ws.Range('E33').Value = '  +6.73%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.65'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.40%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.181'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.45%  '

$ws.Range('E37').Value = '  +0.20%  '

$ws.Range('D38').Value = '3.765.96'
$ws.Range('E38').Value = '  +1.14%  '

$ws.Range('E39').Value = '  +3.83%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.06%  '

$ws.Range('E41').Value = '  +4.20%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000342'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +31.53%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.26'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +13.96%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.978'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.56%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '161.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.43%  '

$ws.Range('B48').Value = 'Arweave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.60%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '49.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.25%  '

$ws.Range('E50').Value = '  +1.34%  '

$ws.Range('E51').Value = '  +2.99%  '
